# Generate Report for Archive
#
# 1. Flip the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2/F2, zh-cn!C2,
#    de-de!C2 - the "Status" cells for the two locales).
# 2. Narrow the status column(s) that held that text: Overview columns
#    E & F ("zh-cn"/"de-de" status columns) and column C ("Status") on
#    the per-locale sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth  = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

# --- Update status text ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZh.Range("C2").Value = $newStatus
$wsDe.Range("C2").Value = $newStatus

# --- Narrow the status columns to match the new, shorter text ---
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZh.Columns.Item(3).ColumnWidth = $newWidth
$wsDe.Columns.Item(3).ColumnWidth = $newWidth
